# Generate Report for Handback
# Fills in the "e1587f2b-c23a-4d35-8bc6-b411e087009c" row (row 7) Latest
# Target/Handback columns on the zh-cn and de-de sheets, which previously had
# no handback info recorded, and flags that the handback was generated from a
# stale source revision.

$wb = $excel.ActiveWorkbook

$mdName   = "e1587f2b-c23a-4d35-8bc6-b411e087009c.md"
$hbUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af175ab8ab9a96338448a2b9a4cdb9e36dbb45d7/e2e/e1587f2b-c23a-4d35-8bc6-b411e087009c.md"
$errMsg   = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bfdc6f65c9f9337b9b533a3af9a7c765f70053e/e2e/e1587f2b-c23a-4d35-8bc6-b411e087009c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af175ab8ab9a96338448a2b9a4cdb9e36dbb45d7/e2e/e1587f2b-c23a-4d35-8bc6-b411e087009c.md."

# ----------------------------------------------------------------------
# zh-cn sheet, row 7 (e1587f2b-...)
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(7, 9).Value = $mdName
$wsZh.Cells.Item(7, 9).Font.Underline = 2
$wsZh.Cells.Item(7, 9).Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(7, 9), $hbUrl, [Type]::Missing, [Type]::Missing, $mdName)

$wsZh.Cells.Item(7, 10).Value = "e1587f2b-c23a-4d35-8bc6-b411e087009c.90128d2ce0692cb3e941289121e4b0e019c68e53.zh-cn.xlf"
$wsZh.Cells.Item(7, 11).Value = "2016-08-15 14:53:03"
$wsZh.Cells.Item(7, 16).Value = $errMsg

# ----------------------------------------------------------------------
# de-de sheet, row 7 (e1587f2b-...)
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(7, 9).Value = $mdName
$wsDe.Cells.Item(7, 9).Font.Underline = 2
$wsDe.Cells.Item(7, 9).Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(7, 9), $hbUrl, [Type]::Missing, [Type]::Missing, $mdName)

$wsDe.Cells.Item(7, 10).Value = "e1587f2b-c23a-4d35-8bc6-b411e087009c.90128d2ce0692cb3e941289121e4b0e019c68e53.de-de.xlf"
$wsDe.Cells.Item(7, 11).Value = "2016-08-15 14:53:14"
$wsDe.Cells.Item(7, 16).Value = $errMsg
